# Fix Training Data Issue (#48)
# The "Date" column (BF) stored dates in a malformed "4-28-2013-14" format
# (day-month mashed with the season string). Correct it to the real ISO
# date "2014-04-28" for every data row (rows 2-31).
#
# The value must stay plain text (matching the source workbook's
# t="inlineStr"/shared-string representation) rather than being
# auto-converted into a date serial number, so the target range is
# pre-formatted as Text before the values are written - mirroring how a
# user would type a date-looking string into a text column in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $current = $cell.Value()
    if ($current -eq "4-28-2013-14") {
        $cell.Value = "2014-04-28"
    }
}
